$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 93; this shifts the existing rows 93:134 down to 94:135
# and the sheet's used range grows to A1:R135 automatically.
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new weekly data point.
$ws.Cells.Item(93, 1).Value = 6
$ws.Cells.Item(93, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(93, 3).Value = "Metropolitana"
$ws.Cells.Item(93, 4).Value = 44553
$ws.Cells.Item(93, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(93, 5).Value = 13
$ws.Cells.Item(93, 6).Value = 100112029
$ws.Cells.Item(93, 7).Value = "Orégano"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 35
$ws.Cells.Item(93, 11).Value = 10000
$ws.Cells.Item(93, 12).Value = 11000
$ws.Cells.Item(93, 13).Value = 10457
$ws.Cells.Item(93, 14).Value = "$/docena de atados"
$ws.Cells.Item(93, 15).Value = "Región Metropolitana"
$ws.Cells.Item(93, 16).Value = 3486
$ws.Cells.Item(93, 17).Value = 3
$ws.Cells.Item(93, 18).Value = "Hortaliza"
